# Update the "想去人数" (interested-count) column F on sheets that list
# exhibition/show entries. Values mirror a refreshed data pull where each
# count ticked up slightly.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 53
$ws1.Range("F4").Value  = 624
$ws1.Range("F5").Value  = 163
$ws1.Range("F6").Value  = 9411
$ws1.Range("F7").Value  = 845
$ws1.Range("F9").Value  = 1202
$ws1.Range("F10").Value = 1144
$ws1.Range("F15").Value = 422
$ws1.Range("F16").Value = 88
$ws1.Range("F17").Value = 252
$ws1.Range("F18").Value = 1280

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# Sheet "全部类型" (all types, aggregate of the above)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 53
$ws4.Range("F4").Value  = 12
$ws4.Range("F5").Value  = 624
$ws4.Range("F6").Value  = 163
$ws4.Range("F7").Value  = 9411
$ws4.Range("F8").Value  = 845
$ws4.Range("F10").Value = 1202
$ws4.Range("F11").Value = 1144
$ws4.Range("F16").Value = 422
$ws4.Range("F17").Value = 88
$ws4.Range("F18").Value = 252
$ws4.Range("F19").Value = 1280
